# Applies the commit "commit into new b": wraps several misspelled /
# grammar-flagged runs with <w:proofErr/> markers (as Word's proofing
# pass would on save) and appends a new run of text after "321".

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1: "Sdf.,kdfskfjsdkfslkfjsl kj" -------------------------
# Split into runs with proofErr wrapping:
#   [spellStart]Sdf[spellEnd] . [gramStart] , [spellStart]kdfskfjsdkfslkfjsl[spellEnd][gramEnd]  [spellStart]kj[spellEnd]
$p1 = $d.Paragraphs(1)
$xml1 = '<w:p ' + $wns + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Sdf</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>,</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>kdfskfjsdkfslkfjsl</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>kj</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
$p1.Range.InsertXML($xml1)

# --- Paragraph 3: "K" + "hkghggjhgjhj" ----------------------------------
$p3 = $d.Paragraphs(3)
$xml3 = '<w:p ' + $wns + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>K</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>hkghggjhgjhj</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
$p3.Range.InsertXML($xml3)

# --- Paragraph 4: "Sdjkshhdsddsad" --------------------------------------
$p4 = $d.Paragraphs(4)
$xml4 = '<w:p ' + $wns + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Sdjkshhdsddsad</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
$p4.Range.InsertXML($xml4)

# --- Paragraph 5: "Sdjkasjdlkajsdkjasd" ---------------------------------
$p5 = $d.Paragraphs(5)
$xml5 = '<w:p ' + $wns + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Sdjkasjdlkajsdkjasd</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
$p5.Range.InsertXML($xml5)

# --- Paragraph 6: "L" + "jsldkjsalkdjaklsdkljaskld" ---------------------
$p6 = $d.Paragraphs(6)
$xml6 = '<w:p ' + $wns + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>L</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>jsldkjsalkdjaklsdkljaskld</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
$p6.Range.InsertXML($xml6)

# --- Paragraph 11: "321" -> append new run ------------------------------
$p11 = $d.Paragraphs(11)
$xml11 = '<w:p ' + $wns + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>321</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>lhhkjhjkkkkkkkkkkkkkkkkkkkkkkk</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$p11.Range.InsertXML($xml11)
